$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.173.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.958.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.74%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.89"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -11.40%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.369"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "56.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0784"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.37%  "
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.842"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "13.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.64%  "
$ws.Range("E15").Value = "  +3.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.244.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.955.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.055.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0846"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.43%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.55%  "
$ws.Range("E26").Value = "  -5.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.123"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.62%  "
$ws.Range("E33").Value = "  -7.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0610"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.47%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.74%  "
$ws.Range("E38").Value = "  -3.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0984"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.43%  "
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("E44").Value = "  -3.98%  "
$ws.Range("E45").Value = "  -6.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.334.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.70%  "
$ws.Range("E50").Value = "  -4.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.135.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.53%  "
